$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 13.48457052857455
$ws.Range("D2").Value = 5.176668469940501
$ws.Range("E2").Value = 14.10757829808751
$ws.Range("F2").Value = 27.15801106511764
$ws.Range("G2").Value = 33.98184631538024
$ws.Range("H2").Value = 15.3679665752588
$ws.Range("I2").Value = 26.06348915289791
$ws.Range("K2").Value = 9.899114111532027
$ws.Range("L2").Value = 9.374449907358573
$ws.Range("N2").Value = 19.62611424727102
$ws.Range("B3").Value = 13.3598171183138
$ws.Range("D3").Value = 5.160803378983115
$ws.Range("E3").Value = 14.13050070204311
$ws.Range("F3").Value = 27.08962549683395
$ws.Range("G3").Value = 33.8286027481911
$ws.Range("H3").Value = 15.39610713573379
$ws.Range("I3").Value = 26.15367265503337
$ws.Range("K3").Value = 9.573744014197953
$ws.Range("L3").Value = 9.33288049383826
$ws.Range("N3").Value = 19.6900717145537
$ws.Range("B4").Value = 13.28533463329183
$ws.Range("D4").Value = 5.150876610175571
$ws.Range("E4").Value = 14.14552958154498
$ws.Range("F4").Value = 27.05539240173879
$ws.Range("G4").Value = 33.74571990956162
$ws.Range("H4").Value = 15.41638499060262
$ws.Range("I4").Value = 26.21341704443844
$ws.Range("K4").Value = 9.37031886247274
$ws.Range("L4").Value = 9.309188732026564
$ws.Range("N4").Value = 19.73114908080441
$ws.Range("B5").Value = 13.25554569527032
$ws.Range("D5").Value = 5.146785479991697
$ws.Range("E5").Value = 14.151894477968
$ws.Range("F5").Value = 27.04340156231003
$ws.Range("G5").Value = 33.71478929362857
$ws.Range("H5").Value = 15.42540142583807
$ws.Range("I5").Value = 26.2388619031735
$ws.Range("K5").Value = 9.286642077663819
$ws.Range("L5").Value = 9.300002362031014
$ws.Range("N5").Value = 19.7483442362025
$ws.Range("B6").Value = 13.25063417224965
$ws.Range("D6").Value = 5.146103401532206
$ws.Range("E6").Value = 14.15296590511346
$ws.Range("F6").Value = 27.04152905889603
$ws.Range("G6").Value = 33.70982577916244
$ws.Range("H6").Value = 15.42694404252248
$ws.Range("I6").Value = 26.24315333179207
$ws.Range("K6").Value = 9.272704768735288
$ws.Range("L6").Value = 9.298505459783838
$ws.Range("N6").Value = 19.75122704739082
$ws.Range("B7").Value = 13.28493056831951
$ws.Range("D7").Value = 5.150821620576981
$ws.Range("E7").Value = 14.14561444626517
$ws.Range("F7").Value = 27.05522274513116
$ws.Range("G7").Value = 33.74529121895809
$ws.Range("H7").Value = 15.41650354205458
$ws.Range("I7").Value = 26.2137557556828
$ws.Range("K7").Value = 9.369193334119561
$ws.Range("L7").Value = 9.309062936134064
$ws.Range("N7").Value = 19.73137913340283
$ws.Range("B8").Value = 13.44113517387515
$ws.Range("D8").Value = 5.171237151039152
$ws.Range("E8").Value = 14.11528415502179
$ws.Range("F8").Value = 27.13282814822223
$ws.Range("G8").Value = 33.92669854735763
$ws.Range("H8").Value = 15.37704623370201
$ws.Range("I8").Value = 26.09367674364303
$ws.Range("K8").Value = 9.787762553624505
$ws.Range("L8").Value = 9.35974123732878
$ws.Range("N8").Value = 19.64779251232003
$ws.Range("B9").Value = 13.7628381313277
$ws.Range("D9").Value = 5.209767252606843
$ws.Range("E9").Value = 14.06335687994873
$ws.Range("F9").Value = 27.34604631086015
$ws.Range("G9").Value = 34.37000779596629
$ws.Range("H9").Value = 15.32351808537059
$ws.Range("I9").Value = 25.89291713265176
$ws.Range("K9").Value = 10.57415045483488
$ws.Range("L9").Value = 9.47329971008355
$ws.Range("N9").Value = 19.49815589828879
$ws.Range("B10").Value = 14.0066379640546
$ws.Range("D10").Value = 5.237115919149676
$ws.Range("E10").Value = 14.02977777689156
$ws.Range("F10").Value = 27.53907979483082
$ws.Range("G10").Value = 34.74694818461435
$ws.Range("H10").Value = 15.29878788078443
$ws.Range("I10").Value = 25.76662349422828
$ws.Range("K10").Value = 11.12434979947011
$ws.Range("L10").Value = 9.564858958218032
$ws.Range("N10").Value = 19.39683130058738
$ws.Range("B11").Value = 14.118744992267
$ws.Range("D11").Value = 5.249340327018202
$ws.Range("E11").Value = 14.01548800906574
$ws.Range("F11").Value = 27.63457349322811
$ws.Range("G11").Value = 34.92901043061406
$ws.Range("H11").Value = 15.29071639967584
$ws.Range("I11").Value = 25.71378274038551
$ws.Range("K11").Value = 11.36749340855789
$ws.Range("L11").Value = 9.6081561525046
$ws.Range("N11").Value = 19.35258722491817
$ws.Range("B12").Value = 14.16133430288896
$ws.Range("D12").Value = 5.253937452745299
$ws.Range("E12").Value = 14.01021807124887
$ws.Range("F12").Value = 27.6718170371582
$ws.Range("G12").Value = 34.99942408599708
$ws.Range("H12").Value = 15.28811743248314
$ws.Range("I12").Value = 25.69443721527226
$ws.Range("K12").Value = 11.45845437772809
$ws.Range("L12").Value = 9.624777593416889
$ws.Range("N12").Value = 19.33609761125181
$ws.Range("B13").Value = 14.15215641830142
$ws.Range("D13").Value = 5.252948815682275
$ws.Range("E13").Value = 14.01134676967759
$ws.Range("F13").Value = 27.66374822366416
$ws.Range("G13").Value = 34.9841947051913
$ws.Range("H13").Value = 15.28865681232303
$ws.Range("I13").Value = 25.69857406590138
$ws.Range("K13").Value = 11.438915021339
$ws.Range("L13").Value = 9.621188008707204
$ws.Range("N13").Value = 19.33963719326711
$ws.Range("B14").Value = 14.12224630376686
$ws.Range("D14").Value = 5.249719181296165
$ws.Range("E14").Value = 14.0150516189537
$ws.Range("F14").Value = 27.63761598990024
$ws.Range("G14").Value = 34.93477424161817
$ws.Range("H14").Value = 15.2904934094894
$ws.Range("I14").Value = 25.71217785616686
$ws.Range("K14").Value = 11.37499956487044
$ws.Range("L14").Value = 9.609519148908936
$ws.Range("N14").Value = 19.35122531732643
$ws.Range("B15").Value = 14.10394224835063
$ws.Range("D15").Value = 5.247736745747849
$ws.Range("E15").Value = 14.01733933316755
$ws.Range("F15").Value = 27.62174946909888
$ws.Range("G15").Value = 34.90469269397575
$ws.Range("H15").Value = 15.29167797279719
$ws.Range("I15").Value = 25.72059709106301
$ws.Range("K15").Value = 11.33570230291248
$ws.Range("L15").Value = 9.602400698815437
$ws.Range("N15").Value = 19.35835780795937
$ws.Range("B16").Value = 13.99933260436552
$ws.Range("D16").Value = 5.236312594197905
$ws.Range("E16").Value = 14.03073144126993
$ws.Range("F16").Value = 27.53299175883084
$ws.Range("G16").Value = 34.73525870340451
$ws.Range("H16").Value = 15.29937937376007
$ws.Range("I16").Value = 25.770169636683
$ws.Range("K16").Value = 11.10830854820466
$ws.Range("L16").Value = 9.562061581911417
$ws.Range("N16").Value = 19.39975985644604
$ws.Range("B17").Value = 13.93544005095161
$ws.Range("D17").Value = 5.229248162989119
$ws.Range("E17").Value = 14.03919917580617
$ws.Range("F17").Value = 27.48049392100769
$ws.Range("G17").Value = 34.63399276935839
$ws.Range("H17").Value = 15.30491839901979
$ws.Range("I17").Value = 25.80176250952343
$ws.Range("K17").Value = 10.96691642973433
$ws.Range("L17").Value = 9.537728814966744
$ws.Range("N17").Value = 19.42563136820299
$ws.Range("B18").Value = 13.8988059605378
$ws.Range("D18").Value = 5.225164548160476
$ws.Range("E18").Value = 14.04416238140966
$ws.Range("F18").Value = 27.4510230776233
$ws.Range("G18").Value = 34.57674743244934
$ws.Range("H18").Value = 15.3084034502668
$ws.Range("I18").Value = 25.82036782452293
$ws.Range("K18").Value = 10.88492334897376
$ws.Range("L18").Value = 9.523888907099604
$ws.Range("N18").Value = 19.44068608587317
$ws.Range("B19").Value = 13.88642317217734
$ws.Range("D19").Value = 5.223778438416189
$ws.Range("E19").Value = 14.04585878731902
$ws.Range("F19").Value = 27.44116984606242
$ws.Range("G19").Value = 34.55753845704581
$ws.Range("H19").Value = 15.30963479060345
$ws.Range("I19").Value = 25.82674175944753
$ws.Range("K19").Value = 10.85704990240276
$ws.Range("L19").Value = 9.519230015575037
$ws.Range("N19").Value = 19.44581329541504
$ws.Range("B20").Value = 13.94222985855721
$ws.Range("D20").Value = 5.230002295034162
$ws.Range("E20").Value = 14.03828817089772
$ws.Range("F20").Value = 27.48600757632205
$ws.Range("G20").Value = 34.64466955665858
$ws.Range("H20").Value = 15.30429779531823
$ws.Range("I20").Value = 25.79835447626334
$ws.Range("K20").Value = 10.98203765037197
$ws.Range("L20").Value = 9.540303048331092
$ws.Range("N20").Value = 19.4228592909379
$ws.Range("B21").Value = 14.13102820190091
$ws.Range("D21").Value = 5.250668678146954
$ws.Range("E21").Value = 14.01395958383261
$ws.Range("F21").Value = 27.645262480133
$ws.Range("G21").Value = 34.94925074975119
$ws.Range("H21").Value = 15.28994153635179
$ws.Range("I21").Value = 25.70816406045313
$ws.Range("K21").Value = 11.39380389798168
$ws.Range("L21").Value = 9.612940539283102
$ws.Range("N21").Value = 19.34781443136961
$ws.Range("B22").Value = 14.25519976937179
$ws.Range("D22").Value = 5.2639884907681
$ws.Range("E22").Value = 13.99888277948689
$ws.Range("F22").Value = 27.75564057177641
$ws.Range("G22").Value = 35.15685766503334
$ws.Range("H22").Value = 15.28322579950185
$ws.Range("I22").Value = 25.65309092623277
$ws.Range("K22").Value = 11.70726560484622
$ws.Range("L22").Value = 9.661723713596023
$ws.Range("N22").Value = 19.30031034340416
$ws.Range("B23").Value = 14.18886742902889
$ws.Range("D23").Value = 5.256896831125035
$ws.Range("E23").Value = 14.00685435955139
$ws.Range("F23").Value = 27.69616157638657
$ws.Range("G23").Value = 35.04529012392241
$ws.Range("H23").Value = 15.2865659838007
$ws.Range("I23").Value = 25.68212989030579
$ws.Range("K23").Value = 11.52916510722073
$ws.Range("L23").Value = 9.635571053118108
$ws.Range("N23").Value = 19.32552347208282
$ws.Range("B24").Value = 13.9391598766998
$ws.Range("D24").Value = 5.229661420942979
$ws.Range("E24").Value = 14.03869974033223
$ws.Range("F24").Value = 27.4835126367328
$ws.Range("G24").Value = 34.63983954382014
$ws.Range("H24").Value = 15.30457743392784
$ws.Range("I24").Value = 25.79989387087682
$ws.Range("K24").Value = 10.97520353547583
$ws.Range("L24").Value = 9.539138771463998
$ws.Range("N24").Value = 19.42411198363831
$ws.Range("B25").Value = 13.67435950841486
$ws.Range("D25").Value = 5.199509270393776
$ws.Range("E25").Value = 14.07659949799086
$ws.Range("F25").Value = 27.28190737053887
$ws.Range("G25").Value = 34.24091023911562
$ws.Range("H25").Value = 15.33543909731818
$ws.Range("I25").Value = 25.94350722724192
$ws.Range("K25").Value = 10.3658185611294
$ws.Range("L25").Value = 9.44111238363897
$ws.Range("N25").Value = 19.53711750372055
